$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hash-map parameters (mod1, mod2) and the first key value.
$ws.Range("B2").Value = 7
$ws.Range("B3").Value = 4
$ws.Range("A6").Value = 6

# Move the active-cell selection from A6 to A7, matching the saved view state.
$ws.Range("A7").Select()
